$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date value for every data row
# (rows 2 through 533). The value changes from serial date 45172
# (2023-09-03) to 45175 (2023-09-06) for all of them.
$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(45175)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 533 }

$ws.Range("C2:C$lastRow").Value = $newDate
